$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 20.91846350703119
$ws.Range("C2").Value = 14.9376749583435
$ws.Range("E2").Value = 16.4550868695369
$ws.Range("F2").Value = 16.86991607391245
$ws.Range("G2").Value = 3.629624543719184
$ws.Range("N2").Value = 18.1777315890429
$ws.Range("B3").Value = 20.01735535115251
$ws.Range("C3").Value = 14.05222205956468
$ws.Range("E3").Value = 15.5133573386915
$ws.Range("F3").Value = 15.89584955866815
$ws.Range("G3").Value = 3.635678044690552
$ws.Range("N3").Value = 18.19586449334508
$ws.Range("B4").Value = 19.45140641590864
$ws.Range("C4").Value = 13.48445925691748
$ws.Range("E4").Value = 14.91211500379413
$ws.Range("F4").Value = 15.26997757108489
$ws.Range("G4").Value = 3.639567762208798
$ws.Range("N4").Value = 18.20887109872193
$ws.Range("B5").Value = 19.21797740009319
$ws.Range("C5").Value = 13.24725046560427
$ws.Range("E5").Value = 14.66158620008865
$ws.Range("F5").Value = 15.00819731993643
$ws.Range("G5").Value = 3.641196585913178
$ws.Range("N5").Value = 18.21463764266039
$ws.Range("B6").Value = 19.17905967776397
$ws.Range("C6").Value = 13.20751662614565
$ws.Range("E6").Value = 14.61966178064999
$ws.Range("F6").Value = 14.96433081551586
$ws.Range("G6").Value = 3.641469700022393
$ws.Range("N6").Value = 18.2156231956182
$ws.Range("B7").Value = 19.44826910759197
$ws.Range("C7").Value = 13.48128350437767
$ws.Range("E7").Value = 14.90875822944741
$ws.Range("F7").Value = 15.26647399323726
$ws.Range("G7").Value = 3.639589551667163
$ws.Range("N7").Value = 18.20894698673534
$ws.Range("B8").Value = 20.61060690338662
$ws.Range("C8").Value = 14.63748730653015
$ws.Range("E8").Value = 16.13530170532846
$ws.Range("F8").Value = 16.53996406344768
$ws.Range("G8").Value = 3.631676080358685
$ws.Range("N8").Value = 18.18359224130453
$ws.Range("B9").Value = 22.77479932307783
$ws.Range("C9").Value = 16.70624258055679
$ws.Range("E9").Value = 18.42966497513565
$ws.Range("F9").Value = 19.0027458068253
$ws.Range("G9").Value = 3.617516841641076
$ws.Range("N9").Value = 18.14894610506398
$ws.Range("B10").Value = 24.2779365995826
$ws.Range("C10").Value = 18.09762428175459
$ws.Range("E10").Value = 20.10457653876821
$ws.Range("F10").Value = 20.67494806633232
$ws.Range("G10").Value = 3.607925082908879
$ws.Range("N10").Value = 18.13297830831772
$ws.Range("B11").Value = 24.94011825281613
$ws.Range("C11").Value = 18.70166202549446
$ws.Range("E11").Value = 20.82582255864779
$ws.Range("F11").Value = 21.3917225636224
$ws.Range("G11").Value = 3.603733835850118
$ws.Range("N11").Value = 18.127838811344
$ws.Range("B12").Value = 25.18757663087733
$ws.Range("C12").Value = 18.92617935517882
$ws.Range("E12").Value = 21.09313844167677
$ws.Range("F12").Value = 21.6568656903351
$ws.Range("G12").Value = 3.602171161565135
$ws.Range("N12").Value = 18.12620345314721
$ws.Range("B13").Value = 25.13443130556133
$ws.Range("C13").Value = 18.87801396228673
$ws.Range("E13").Value = 21.0358246087094
$ws.Range("F13").Value = 21.60004134736749
$ws.Range("G13").Value = 3.602506628169308
$ws.Range("N13").Value = 18.1265417450479
$ws.Range("B14").Value = 24.96054382624261
$ws.Range("C14").Value = 18.72021804320799
$ws.Range("E14").Value = 20.8479308641777
$ws.Range("F14").Value = 21.41366180504534
$ws.Range("G14").Value = 3.60360478513743
$ws.Range("N14").Value = 18.12769800840678
$ws.Range("B15").Value = 24.85359845468544
$ws.Range("C15").Value = 18.62301255034317
$ws.Range("E15").Value = 20.7320859737739
$ws.Range("F15").Value = 21.29868154950791
$ws.Range("G15").Value = 3.604280614834169
$ws.Range("N15").Value = 18.1284468965429
$ws.Range("B16").Value = 24.23420992984627
$ws.Range("C16").Value = 18.05756161916133
$ws.Range("E16").Value = 20.05662671023068
$ws.Range("F16").Value = 20.62722412089977
$ws.Range("G16").Value = 3.608202429273535
$ws.Range("N16").Value = 18.13335739232222
$ws.Range("B17").Value = 23.84855805986709
$ws.Range("C17").Value = 17.70322218638001
$ws.Range("E17").Value = 19.6318675221266
$ws.Range("F17").Value = 20.20408069617459
$ws.Range("G17").Value = 3.61065220949227
$ws.Range("N17").Value = 18.13691772699663
$ws.Range("B18").Value = 23.62471720574018
$ws.Range("C18").Value = 17.49669696214199
$ws.Range("E18").Value = 19.38372558959297
$ws.Range("F18").Value = 19.95656407809808
$ws.Range("G18").Value = 3.6120774737588
$ws.Range("N18").Value = 18.13916516482735
$ws.Range("B19").Value = 23.54858714657355
$ws.Range("C19").Value = 17.42630631587686
$ws.Range("E19").Value = 19.2990495214094
$ws.Range("F19").Value = 19.87204792380562
$ws.Range("G19").Value = 3.612562837442619
$ws.Range("N19").Value = 18.13996023161808
$ws.Range("B20").Value = 23.88982230921167
$ws.Range("C20").Value = 17.74122408391776
$ws.Range("E20").Value = 19.67748009911441
$ws.Range("F20").Value = 20.2495528364879
$ws.Range("G20").Value = 3.610389750031298
$ws.Range("N20").Value = 18.13651802381442
$ws.Range("B21").Value = 25.01170959201154
$ws.Range("C21").Value = 18.76668143143275
$ws.Range("E21").Value = 20.90327697741647
$ws.Range("F21").Value = 21.46857628470567
$ws.Range("G21").Value = 3.603281568279502
$ws.Range("N21").Value = 18.12734990605808
$ws.Range("B22").Value = 25.72565179819505
$ws.Range("C22").Value = 19.41989667140476
$ws.Range("E22").Value = 21.67059516751293
$ws.Range("F22").Value = 22.22866616901555
$ws.Range("G22").Value = 3.598778396662258
$ws.Range("N22").Value = 18.12317246442178
$ws.Range("B23").Value = 25.34642615550571
$ws.Range("C23").Value = 19.06997500409409
$ws.Range("E23").Value = 21.26414168747459
$ws.Range("F23").Value = 21.82633154475864
$ws.Range("G23").Value = 3.601168886836671
$ws.Range("N23").Value = 18.12523422879307
$ws.Range("B24").Value = 23.87117334464579
$ws.Range("C24").Value = 17.72405217609973
$ws.Range("E24").Value = 19.65687093561468
$ws.Range("F24").Value = 20.22900810905294
$ws.Range("G24").Value = 3.610508355438539
$ws.Range("N24").Value = 18.13669810534926
$ws.Range("B25").Value = 22.20353165852204
$ws.Range("C25").Value = 16.16873454922522
$ws.Range("E25").Value = 17.77649002815375
$ws.Range("F25").Value = 18.34778573295697
$ws.Range("G25").Value = 3.62120359135831
$ws.Range("N25").Value = 18.15667447208433
